$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "Type" column for each hero row with "h" first so that the
# shared-string table gets "h" appended before "Type"/"Mode"
$ws.Range("E2").Value = "h"
$ws.Range("E3").Value = "h"
$ws.Range("E4").Value = "h"
$ws.Range("E5").Value = "h"

# Add new headers for Type and Mode columns
$ws.Range("E1").Value = "Type"
$ws.Range("F1").Value = "Mode"

# Update the selection to match the target state
$ws.Range("E5").Select()
